$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Permits Filed for 35-34 41st Street in Astoria, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2025/10/permits-filed-for-35-34-41st-street-in-astoria-queens.html"
$ws.Range("C2").Value = 'Permits have been filed for a 12-story mixed-use building at 35-34 41st Street in <a href="https://newyorkyimby.com/neighborhoods/astoria">Astoria</a>, Queens. Located between 35th and 36th Avenues, the lot is near the Steinway Street subway station, served by the E, F, M, and R trains. Stephen Ohnemus of Domain 41st Street Site B LLC is listed as the owner behind the applications.'
$ws.Range("D2").Value = "2025-10-11T11:00:47+00:00"
$ws.Range("E2").Value = "Sat, 11 Oct 2025 11:00:47 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Astoria"
$ws.Cells.Item(2, 8).Value = "'"
$ws.Cells.Item(2, 8).Style = "Normal"
